$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.804.82'
$ws.Range('E2').Value = '  -1.90%  '
$ws.Range('D3').Value = '2.671.95'
$ws.Range('E3').Value = '  -2.34%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '551.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.75'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.591'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.94%  '
$ws.Range('E9').Value = '  -3.27%  '
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('E11').Value = '  -3.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.26'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -7.93%  '
$ws.Range('D13').Value = '3.146.88'
$ws.Range('E13').Value = '  -2.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.12'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.77%  '
$ws.Range('D15').Value = '62.749.71'
$ws.Range('E15').Value = '  -1.39%  '
$ws.Range('E16').Value = '  -2.54%  '
$ws.Range('D17').Value = '2.675.77'
$ws.Range('E17').Value = '  -2.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.86'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.59'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '344.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.28'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('E23').Value = '  -3.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.07%  '
$ws.Range('E25').Value = '  -1.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.995'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.57%  '
$ws.Range('E27').Value = '  -3.05%  '
$ws.Range('D28').Value = '0.0₃0857'
$ws.Range('E28').Value = '  -5.72%  '
$ws.Range('E29').Value = '  +1.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('E31').Value = '  -1.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '166.35'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.49'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.40%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.85'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.79%  '
$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.77'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '348.75'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.958'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.03%  '
$ws.Range('E40').Value = '  -0.78%  '
$ws.Range('E41').Value = '  -3.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '38.23'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.29'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.08%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.76'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.94%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0562'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.615'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.41%  '
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '11.03'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0972'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.07%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '128.91'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.20%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0241'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.83%  '
